$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# FilesTab query (B4): removed the "File Type" and "Breed" columns
# from the RETURN clause (ICDC Breed script correction).
$fileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $fileQuery

# Row 4 height shrinks because the corrected query text has two fewer lines.
$ws.Rows.Item(4).RowHeight = 217.5

# Scroll/select so row 4 (FilesTab) is in view, matching the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B4").Select()
